# This presentation ships with two theme parts:
#   ppt/theme/theme1.xml  -> currently the "Office Theme" colour scheme,
#                             only used by the notes master.
#   ppt/theme/theme2.xml  -> currently the "Integral" / "Red Violet" colour
#                             scheme, used by the slide master (and so by
#                             every slide).
#
# The authored edit swaps the two themes' contents: theme1.xml ends up
# holding the Integral/Red Violet colours and theme2.xml ends up holding
# the Office colours. Concretely this means the (identical) font/format
# schemes stay put, but every one of the 12 theme colour slots is swapped
# between the slide-master theme and the notes-master theme.
#
# We reach the slide-master theme through Slide.ThemeColorScheme and the
# notes-master theme through Slide.NotesPage.ThemeColorScheme, and push
# the 12 colours (Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink,
# FollowedHyperlink - MsoThemeColorSchemeIndex order 1..12) across.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Office theme colours (were in theme1.xml) -> move onto the slide master
# theme (theme2.xml), replacing the old Red Violet / Integral colours.
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

# Red Violet / Integral theme colours (were in theme2.xml) -> move onto the
# notes master theme (theme1.xml), replacing the old Office colours.
$redVioletColors = @(0, 16777215, 5326149, 14473688, 9514467, 13381832, 14460494, 15168839, 14774665, 7555029, 2465643, 9211020)

$slideScheme = $s.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $slideScheme.Item($i).RGB = $officeColors[$i - 1]
}

$notesScheme = $s.NotesPage.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Item($i).RGB = $redVioletColors[$i - 1]
}
